$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Mario Rossi) ---
$ws.Range("C2").Value = "1/1/1980"
$ws.Range("E2").Value = "Coniugato"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = "Laurea magistrale in Economia"
$ws.Range("J2").Value = "Direttore di banca"

# --- Row 3 (Filippo Bianchi) ---
$ws.Range("C3").Value = "1/1/1977"
$ws.Range("E3").Value = "Celibe"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "Licenza media"
$ws.Range("J3").Value = "Barista"

# --- Row 4 (Giuseppina Cafueri) ---
$ws.Range("C4").Value = "1/1/1965"
$ws.Range("E4").Value = "Coniugata"
$ws.Range("I4").Value = "Laurea magistrale in Giurisprudenza"
$ws.Range("J4").Value = "Avvocato"

# --- Row 5 (was Vincenzo Verdi -> Giovanni Mele) ---
$ws.Range("A5").Value = "Giovanni"
$ws.Range("B5").Value = "Mele"
$ws.Range("C5").Value = "1/1/1997"
$ws.Range("I5").Value = "Laurea triennale in Informatica"
$ws.Range("J5").Value = "Consulente "
